# Auto-generated edit script: applies the Jenova_Profits.xlsx commit
# (per-worksheet currentAveragePrice / LevePrice* / LeveProfit* refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 37.1
$ws.Range("I11").Value = 37.1
$ws.Range("K11").Value = 37.1
$ws.Range("M11").Value = 102.9
$ws.Range("H41").Value = 20002908
$ws.Range("I41").Value = 478.23077
$ws.Range("J41").Value = 41672210
$ws.Range("K41").Value = 478.23077
$ws.Range("L41").Value = 41672210
$ws.Range("M41").Value = -38.23077000000001
$ws.Range("N41").Value = -41673090
$ws.Range("H62").Value = 12506165
$ws.Range("I62").Value = 31253444
$ws.Range("K62").Value = 31253444
$ws.Range("M62").Value = -31252820
$ws.Range("H65").Value = 12506165
$ws.Range("I65").Value = 31253444
$ws.Range("K65").Value = 156267220
$ws.Range("M65").Value = -156264100
$ws.Range("H70").Value = 201400
$ws.Range("I70").Value = 1500
$ws.Range("K70").Value = 4500
$ws.Range("M70").Value = -4230
$ws.Range("H73").Value = 201400
$ws.Range("I73").Value = 1500
$ws.Range("K73").Value = 4500
$ws.Range("M73").Value = -3564
$ws.Range("H98").Value = 2248.919
$ws.Range("I98").Value = 1578.6129
$ws.Range("K98").Value = 1578.6129
$ws.Range("M98").Value = -80.61290000000008
$ws.Range("H111").Value = 59494.61
$ws.Range("I111").Value = 103852.1
$ws.Range("K111").Value = 311556.3
$ws.Range("M111").Value = -308489.3
$ws.Range("H113").Value = 4256.2856
$ws.Range("I113").Value = 3632.6667
$ws.Range("J113").Value = 4724
$ws.Range("K113").Value = 3632.6667
$ws.Range("L113").Value = 4724
$ws.Range("M113").Value = -378.6667000000002
$ws.Range("N113").Value = -11232
$ws.Range("H122").Value = 2248.919
$ws.Range("I122").Value = 1578.6129
$ws.Range("K122").Value = 4735.8387
$ws.Range("M122").Value = -2285.8387
$ws.Range("H138").Value = 6812
$ws.Range("J138").Value = 7759.756
$ws.Range("L138").Value = 23279.268
$ws.Range("N138").Value = -33559.268
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1113.4814
$ws.Range("I2").Value = 1243.8182
$ws.Range("J2").Value = 540
$ws.Range("K2").Value = 1243.8182
$ws.Range("L2").Value = 540
$ws.Range("M2").Value = -1130.8182
$ws.Range("N2").Value = -766
$ws.Range("H32").Value = 3557.6765
$ws.Range("I32").Value = 3557.6765
$ws.Range("K32").Value = 3557.6765
$ws.Range("M32").Value = -3270.6765
$ws.Range("H116").Value = 1113.4814
$ws.Range("I116").Value = 1243.8182
$ws.Range("J116").Value = 540
$ws.Range("K116").Value = 1243.8182
$ws.Range("L116").Value = 540
$ws.Range("M116").Value = 1050.1818
$ws.Range("N116").Value = -5128
$ws.Range("H118").Value = 126000
$ws.Range("J118").Value = 126000
$ws.Range("L118").Value = 126000
$ws.Range("N118").Value = -129314
$ws.Range("H132").Value = 2002.8096
$ws.Range("I132").Value = 1415.3784
$ws.Range("K132").Value = 4246.135200000001
$ws.Range("M132").Value = -1716.135200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1113.4814
$ws.Range("I3").Value = 1243.8182
$ws.Range("J3").Value = 540
$ws.Range("K3").Value = 1243.8182
$ws.Range("L3").Value = 540
$ws.Range("M3").Value = -1129.8182
$ws.Range("N3").Value = -768
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H20").Value = 4338.1665
$ws.Range("I20").Value = 3671.6667
$ws.Range("K20").Value = 3671.6667
$ws.Range("M20").Value = -3424.6667
$ws.Range("H134").Value = 28790.025
$ws.Range("I134").Value = 1750.6562
$ws.Range("J134").Value = 173000
$ws.Range("K134").Value = 5251.9686
$ws.Range("L134").Value = 519000
$ws.Range("M134").Value = -2716.9686
$ws.Range("N134").Value = -524070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6006
$ws.Range("I58").Value = 5402.6665
$ws.Range("K58").Value = 5402.6665
$ws.Range("M58").Value = -5199.6665
$ws.Range("H62").Value = 4232.3335
$ws.Range("I62").Value = 2499.5
$ws.Range("K62").Value = 2499.5
$ws.Range("M62").Value = -1875.5
$ws.Range("H65").Value = 4232.3335
$ws.Range("I65").Value = 2499.5
$ws.Range("K65").Value = 12497.5
$ws.Range("M65").Value = -9377.5
$ws.Range("H99").Value = 6284.1763
$ws.Range("J99").Value = 7299.6
$ws.Range("L99").Value = 7299.6
$ws.Range("N99").Value = -10295.6
$ws.Range("H107").Value = 431.81818
$ws.Range("I107").Value = 375.1
$ws.Range("K107").Value = 375.1
$ws.Range("M107").Value = 1544.9
$ws.Range("H126").Value = 6284.1763
$ws.Range("J126").Value = 7299.6
$ws.Range("L126").Value = 21898.8
$ws.Range("N126").Value = -26838.8
$ws.Range("H136").Value = 6006
$ws.Range("I136").Value = 5402.6665
$ws.Range("K136").Value = 16207.9995
$ws.Range("M136").Value = -13657.9995
$ws.Range("H141").Value = 213343.67
$ws.Range("J141").Value = 224974.62
$ws.Range("L141").Value = 224974.62
$ws.Range("N141").Value = -235334.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2580
$ws.Range("I75").Value = 1495
$ws.Range("J75").Value = 3122.5
$ws.Range("K75").Value = 4485
$ws.Range("L75").Value = 9367.5
$ws.Range("M75").Value = -3487
$ws.Range("N75").Value = -11363.5
$ws.Range("H78").Value = 2580
$ws.Range("I78").Value = 1495
$ws.Range("J78").Value = 3122.5
$ws.Range("K78").Value = 13455
$ws.Range("L78").Value = 28102.5
$ws.Range("M78").Value = -8463
$ws.Range("N78").Value = -38086.5
$ws.Range("H122").Value = 48428.332
$ws.Range("I122").Value = 695.875
$ws.Range("J122").Value = 77802.16
$ws.Range("K122").Value = 6262.875
$ws.Range("L122").Value = 700219.4400000001
$ws.Range("M122").Value = -3812.875
$ws.Range("N122").Value = -705119.4400000001
$ws.Range("H129").Value = 18519790
$ws.Range("J129").Value = 1829
$ws.Range("L129").Value = 5487
$ws.Range("N129").Value = -15487
$ws.Range("H131").Value = 2795.6035
$ws.Range("J131").Value = 2962.2263
$ws.Range("L131").Value = 8886.678899999999
$ws.Range("N131").Value = -18966.6789
$ws.Range("H138").Value = 50003516
$ws.Range("I138").Value = 7030
$ws.Range("J138").Value = 100000000
$ws.Range("K138").Value = 21090
$ws.Range("L138").Value = 300000000
$ws.Range("M138").Value = -15950
$ws.Range("N138").Value = -300010280
$ws.Range("H139").Value = 5430.3184
$ws.Range("I139").Value = 4676.2354
$ws.Range("K139").Value = 14028.7062
$ws.Range("M139").Value = -8888.706199999999
$ws.Range("H140").Value = 1218.875
$ws.Range("I140").Value = 1218.875
$ws.Range("K140").Value = 3656.625
$ws.Range("M140").Value = 1523.375
$ws.Range("H141").Value = 7384.8887
$ws.Range("I141").Value = 7384.8887
$ws.Range("K141").Value = 22154.6661
$ws.Range("M141").Value = -16974.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 36711.125
$ws.Range("J43").Value = 39949.832
$ws.Range("L43").Value = 39949.832
$ws.Range("N43").Value = -40251.832
$ws.Range("H46").Value = 40059.6
$ws.Range("J46").Value = 48824.5
$ws.Range("L46").Value = 48824.5
$ws.Range("N46").Value = -49136.5
$ws.Range("H132").Value = 617797.75
$ws.Range("J132").Value = 96501.17999999999
$ws.Range("L132").Value = 289503.54
$ws.Range("N132").Value = -294563.54
$ws.Range("H137").Value = 49999
$ws.Range("I137").Value = 49999
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 49999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -44899
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 267112.75
$ws.Range("I7").Value = 3877.1738
$ws.Range("K7").Value = 3877.1738
$ws.Range("M7").Value = -3765.1738
$ws.Range("H126").Value = 267112.75
$ws.Range("I126").Value = 3877.1738
$ws.Range("K126").Value = 11631.5214
$ws.Range("M126").Value = -9161.5214
$ws.Range("H136").Value = 840031.3
$ws.Range("I136").Value = 1256209.2
$ws.Range("K136").Value = 3768627.6
$ws.Range("M136").Value = -3766077.6
$ws.Range("H138").Value = 79970
$ws.Range("J138").Value = 79970
$ws.Range("L138").Value = 79970
$ws.Range("N138").Value = -90250

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1895.2354
$ws.Range("I126").Value = 523.5
$ws.Range("K126").Value = 1570.5
$ws.Range("M126").Value = 899.5
